# Update stats for 2025-07
# Adds three new year-over-year columns (F: yoy_schools, G: yoy_authorities,
# H: yoy_users) and refreshes the latest (July 2025) row's user totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (F1:H1), matching the bold/centered/bordered style
# already used by the other header cells (copy format from E1, then set text) ---
$headerSrc = $ws.Range("E1")
$headerSrc.Copy($ws.Range("F1"))
$headerSrc.Copy($ws.Range("G1"))
$headerSrc.Copy($ws.Range("H1"))

$ws.Range("F1").Value = "yoy_schools"
$ws.Range("G1").Value = "yoy_authorities"
$ws.Range("H1").Value = "yoy_users"

# --- Rows 2-13: the new columns don't have a YoY comparator yet (less than
# 12 months of history), so touch the cells (no-op border assignment) just
# to materialize them as present-but-blank, matching the source layout
# without introducing any new value, type, or style. ---
foreach ($row in 2..13) {
    $ws.Cells.Item($row, 6).Borders.LineStyle = -4142
    $ws.Cells.Item($row, 7).Borders.LineStyle = -4142
    $ws.Cells.Item($row, 8).Borders.LineStyle = -4142
}

# --- Year-over-year values for rows 14-20 (month r vs month r-12) ---
$yoy = @{
    14 = @(2.940153096729303, 5.961754780652417, 23.28061250163025)
    15 = @(3.658536585365857, 6.764374295377684, 24.85887932178075)
    16 = @(4.054289194362282, 6.877113866967299, 25.07756835683654)
    17 = @(5.86376404494382,  3.205128205128216, 18.87096770378025)
    18 = @(6.092436974789917, 3.311965811965822, 18.41667687390272)
    19 = @(6.339254615116685, 3.201707577374591, 24.62859203576528)
    20 = @(6.184142338918641, 3.503184713375807, 25.58277891171774)
}

foreach ($row in $yoy.Keys | Sort-Object) {
    $vals = $yoy[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}

# --- Refresh the July 2025 row (row 20): users (D) and users_per_school (E) ---
$ws.Range("D20").Value = 5551604
$ws.Range("E20").Value = 903.1403936879778

Write-Output "Applied 2025-07 stats update"
